# Adding support for multi-task studies
# Adds a "Comments" column (J) to Sheet1, with a note on row 6 explaining
# a bug found while testing multi-task study support.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J1").Value = "Comments"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J6").Value = "Not writing the correct session/task config record"

$ws.Columns.Item(10).ColumnWidth = 43.6666666666666666
